# Insert a new weekly price record for Coliflor (Macroferia Regional de Talca)
# at row 149, shifting the existing rows 149-248 down to 150-249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 149..248 down by one (Excel copies formatting/number formats
# from the row above automatically, which keeps column D's date format).
$ws.Rows.Item(149).Insert()

# Values for the newly inserted row 149.
$newRowValues = @(
    5,
    "Macroferia Regional de Talca",
    "Maule",
    44719,
    7,
    100112008,
    "Coliflor",
    "Sin especificar",
    "Primera",
    2000,
    1200,
    1200,
    1200,
    "`$/unidad",
    "Región del Maule",
    1200,
    1,
    "Hortaliza"
)

for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item(149, $i + 1).Value = $newRowValues[$i]
}
